$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.637007905307379
$ws.Range("C2").Value = 2.136839981357869
$ws.Range("D2").Value = -0.104899421407004

$ws.Range("B3").Value = 1.524080962532697
$ws.Range("C3").Value = 2.25111111111111
$ws.Range("D3").Value = 0.1698591361477161

$ws.Range("B4").Value = 0.9605357173497543
$ws.Range("C4").Value = 0.7776212832551026
$ws.Range("D4").Value = 0.8690294310658194

$ws.Range("B5").Value = 1.636994083545595
$ws.Range("C5").Value = 2.147574700561443
$ws.Range("D5").Value = -0.1048621058555421
